# Clean up the author-list strings in column E (rows 2-17): each entry is a
# "[Name%Name%email%n, ...]" blob where list items are separated by a comma
# followed by a run of spaces used as pretty-printed indentation. The source
# data regenerated these blobs with one extra space of indentation per
# separator, so every run of 2+ spaces gets exactly one more space appended.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 17; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $old = $cell.Value2
    $new = $old -replace '( {2,})', '$1 '
    $cell.Value = $new
}
